$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bbf1")

# Row 26 (BBF4 table - min row)
$ws.Range("C26").Value = 45.0284
$ws.Range("D26").Value = 46.5566
$ws.Range("E26").Value = 45.3208
$ws.Range("F26").Value = 0.36001
$ws.Range("G26").Value = 43.9823
$ws.Range("H26").Value = 44.0193
$ws.Range("I26").Value = 43.9964
$ws.Range("J26").Value = 0.012371

# Row 27 (BBF4 table - 15100 row)
$ws.Range("C27").Value = 46.6346
$ws.Range("D27").Value = 55.7886
$ws.Range("E27").Value = 49.9422
$ws.Range("F27").Value = 2.7921
$ws.Range("G27").Value = 44.1291
$ws.Range("H27").Value = 44.5433
$ws.Range("I27").Value = 44.3026
$ws.Range("J27").Value = 0.11315

# Row 28 (BBF4 table - mean row)
$ws.Range("C28").Value = 44.1427
$ws.Range("D28").Value = 45.0376
$ws.Range("E28").Value = 44.7143
$ws.Range("F28").Value = 0.27313
$ws.Range("G28").Value = 43.9823
$ws.Range("H28").Value = 44.0289
$ws.Range("I28").Value = 43.9926
$ws.Range("J28").Value = 0.012536

# Row 33 (BBF5 table - min row)
$ws.Range("C33").Value = 68.3066
$ws.Range("D33").Value = 68.7698
$ws.Range("E33").Value = 68.6043
$ws.Range("F33").Value = 0.10064
$ws.Range("G33").Value = 68.2853
$ws.Range("H33").Value = 68.7565
$ws.Range("I33").Value = 68.6155
$ws.Range("J33").Value = 0.10908

# Row 34 (BBF5 table - 15100 row)
$ws.Range("C34").Value = 67.4522
$ws.Range("D34").Value = 68.5608
$ws.Range("E34").Value = 68.334
$ws.Range("F34").Value = 0.28681
$ws.Range("G34").Value = 50.941
$ws.Range("H34").Value = 68.5782
$ws.Range("I34").Value = 64.7729
$ws.Range("J34").Value = 6.2686

# Row 35 (BBF5 table - mean row)
$ws.Range("C35").Value = 67.1249
$ws.Range("D35").Value = 67.3662
$ws.Range("E35").Value = 67.2137
$ws.Range("F35").Value = 0.074606
$ws.Range("G35").Value = 67.1236
$ws.Range("H35").Value = 67.1239
$ws.Range("I35").Value = 67.1238
$ws.Range("J35").Value = 0.000079386000000000002

# New styled (empty) cell Q41 matching the Q5/Q13 numeric-format style
$ws.Range("Q41").NumberFormat = "0.00E+00"

# Restore the authored selection
$ws.Range("R24").Select()
